$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet back to "Sheet1" (was "Investor")
$ws.Name = "Sheet1"

# Revert Investor Code values in column A (rows 2-141) to the old/legacy codes
$ws.Range("A2").Value = "JoSmAEInSBIIPaLP1"
$ws.Range("A3").Value = "JoSmAEMePaIVLP"
$ws.Range("A4").Value = "JoSmAEGrEqPaLP"
$ws.Range("A5").Value = "JoSmAEMiMaDeIIGPLP"
$ws.Range("A6").Value = "JoSmAEInSBIVPaLP"
$ws.Range("A7").Value = "JoSmAEInPaEuL."
$ws.Range("A8").Value = "JoSmAEInSBIIPaLP2"
$ws.Range("A9").Value = "JoSmAEMePaIILP2"
$ws.Range("A10").Value = "JoSmAEInPuEqFuLP"
$ws.Range("A11").Value = "JoSmAEInPaVILP3"
$ws.Range("A12").Value = "JoSmAEInPaVILP2"
$ws.Range("A13").Value = "JoSmAEMiMaDeIVGPLP"
$ws.Range("A14").Value = "JoSmAEInPaVLP"
$ws.Range("A15").Value = "JoSmAEInPaEFIILP"
$ws.Range("A16").Value = "JoSmAEPaEXCFLP"
$ws.Range("A17").Value = "JoSmAEMePaIILP1"
$ws.Range("A18").Value = "JoSmAEGrFuGPLL"
$ws.Range("A19").Value = "JoSmAEImPaLP"
$ws.Range("A20").Value = "JoSmAEInPaVILP1"
$ws.Range("A21").Value = "EmJoAEInPaVILP3"
$ws.Range("A22").Value = "EmJoAEInPaVILP2"
$ws.Range("A23").Value = "EmJoAEMiMaDeIIGPLP"
$ws.Range("A24").Value = "EmJoAEInSBIVPaLP"
$ws.Range("A25").Value = "EmJoAEMiMaDeIVGPLP"
$ws.Range("A26").Value = "EmJoAEInPaEuL."
$ws.Range("A27").Value = "EmJoAEInSBIIPaLP1"
$ws.Range("A28").Value = "EmJoAEInPaVLP"
$ws.Range("A29").Value = "EmJoAEInPaEFIILP"
$ws.Range("A30").Value = "EmJoAEInSBIIPaLP2"
$ws.Range("A31").Value = "EmJoAEMePaIVLP"
$ws.Range("A32").Value = "EmJoAEMePaIILP2"
$ws.Range("A33").Value = "EmJoAEInPuEqFuLP"
$ws.Range("A34").Value = "EmJoAEPaEXCFLP"
$ws.Range("A35").Value = "EmJoAEImPaLP"
$ws.Range("A36").Value = "EmJoAEInPaVILP1"
$ws.Range("A37").Value = "EmJoAEGrEqPaLP"
$ws.Range("A38").Value = "EmJoAEMePaIILP1"
$ws.Range("A39").Value = "MiWiAEInPaVILP3"
$ws.Range("A40").Value = "MiWiAEInPaVILP2"
$ws.Range("A41").Value = "MiWiAEMiMaDeIIGPLP"
$ws.Range("A42").Value = "MiWiAEInSBIVPaLP"
$ws.Range("A43").Value = "MiWiAEMiMaDeIVGPLP"
$ws.Range("A44").Value = "MiWiAEInPaVLP"
$ws.Range("A45").Value = "MiWiAEInPaEFIILP"
$ws.Range("A46").Value = "MiWiAEInSBIIPaLP2"
$ws.Range("A47").Value = "MiWiAEMePaIVLP"
$ws.Range("A48").Value = "MiWiAEMePaIILP2"
$ws.Range("A49").Value = "MiWiAEInPuEqFuLP"
$ws.Range("A50").Value = "MiWiAEPaEXCFLP"
$ws.Range("A51").Value = "MiWiAEImPaLP"
$ws.Range("A52").Value = "MiWiAEInPaVILP1"
$ws.Range("A53").Value = "MiWiAEGrEqPaLP"
$ws.Range("A54").Value = "MiWiAEMePaIILP1"
$ws.Range("A55").Value = "SoBrAEInPaVILP3"
$ws.Range("A56").Value = "SoBrAEInPaVILP2"
$ws.Range("A57").Value = "SoBrAEMiMaDeIIGPLP"
$ws.Range("A58").Value = "SoBrAEInSBIVPaLP"
$ws.Range("A59").Value = "SoBrAEMiMaDeIVGPLP"
$ws.Range("A60").Value = "SoBrAEInSBIIPaLP1"
$ws.Range("A61").Value = "SoBrAEInPaVLP"
$ws.Range("A62").Value = "SoBrAEInPaEFIILP"
$ws.Range("A63").Value = "SoBrAEInSBIIPaLP2"
$ws.Range("A64").Value = "SoBrAEMePaIVLP"
$ws.Range("A65").Value = "SoBrAEMePaIILP2"
$ws.Range("A66").Value = "SoBrAEInPuEqFuLP"
$ws.Range("A67").Value = "SoBrAEPaEXCFLP"
$ws.Range("A68").Value = "SoBrAEImPaLP"
$ws.Range("A69").Value = "SoBrAEInPaVILP1"
$ws.Range("A70").Value = "SoBrAEGrEqPaLP"
$ws.Range("A71").Value = "SoBrAEMePaIILP1"
$ws.Range("A72").Value = "DaDaAEInPaVILP3"
$ws.Range("A73").Value = "DaDaAEInPaVILP2"
$ws.Range("A74").Value = "DaDaAEInSBIVPaLP"
$ws.Range("A75").Value = "DaDaAEMiMaDeIVGPLP"
$ws.Range("A76").Value = "DaDaAEInPaEuL."
$ws.Range("A77").Value = "DaDaAEInPaVLP"
$ws.Range("A78").Value = "DaDaAEInPaEFIILP"
$ws.Range("A79").Value = "DaDaAEInSBIIPaLP2"
$ws.Range("A80").Value = "DaDaAEMePaIVLP"
$ws.Range("A81").Value = "DaDaAEMePaIILP2"
$ws.Range("A82").Value = "DaDaAEPaEXCFLP"
$ws.Range("A83").Value = "DaDaAEImPaLP"
$ws.Range("A84").Value = "DaDaAEInPaVILP1"
$ws.Range("A85").Value = "DaDaAEGrEqPaLP"
$ws.Range("A86").Value = "DaDaAEMePaIILP1"
$ws.Range("A87").Value = "OlMiAEInPaVILP3"
$ws.Range("A88").Value = "OlMiAEInPaVILP2"
$ws.Range("A89").Value = "OlMiAEInSBIVPaLP"
$ws.Range("A90").Value = "OlMiAEInSBIIPaLP1"
$ws.Range("A91").Value = "OlMiAEInPaEuL."
$ws.Range("A92").Value = "OlMiAEInPaVLP"
$ws.Range("A93").Value = "OlMiAEInPaEFIILP"
$ws.Range("A94").Value = "OlMiAEInSBIIPaLP2"
$ws.Range("A95").Value = "OlMiAEMePaIVLP"
$ws.Range("A96").Value = "OlMiAEPaEXCFLP"
$ws.Range("A97").Value = "OlMiAEImPaLP"
$ws.Range("A98").Value = "OlMiAEInPaVILP1"
$ws.Range("A99").Value = "OlMiAEGrEqPaLP"
$ws.Range("A100").Value = "JaWiAEInPaVILP2"
$ws.Range("A101").Value = "JaWiAEMiMaDeIIGPLP"
$ws.Range("A102").Value = "JaWiAEInSBIVPaLP"
$ws.Range("A103").Value = "JaWiAEMiMaDeIVGPLP"
$ws.Range("A104").Value = "JaWiAEInSBIIPaLP1"
$ws.Range("A105").Value = "JaWiAEInPaVLP"
$ws.Range("A106").Value = "JaWiAEInPaEFIILP"
$ws.Range("A107").Value = "JaWiAEInSBIIPaLP2"
$ws.Range("A108").Value = "JaWiAEMePaIVLP"
$ws.Range("A109").Value = "JaWiAEMePaIILP2"
$ws.Range("A110").Value = "JaWiAEInPaVILP1"
$ws.Range("A111").Value = "JaWiAEMePaIILP1"
$ws.Range("A112").Value = "IsMoAEInPaVILP3"
$ws.Range("A113").Value = "IsMoAEInPaVILP2"
$ws.Range("A114").Value = "IsMoAEInSBIVPaLP"
$ws.Range("A115").Value = "IsMoAEMiMaDeIVGPLP"
$ws.Range("A116").Value = "IsMoAEInPaVLP"
$ws.Range("A117").Value = "IsMoAEInPaEFIILP"
$ws.Range("A118").Value = "IsMoAEMePaIVLP"
$ws.Range("A119").Value = "IsMoAEInPuEqFuLP"
$ws.Range("A120").Value = "IsMoAEPaEXCFLP"
$ws.Range("A121").Value = "IsMoAEImPaLP"
$ws.Range("A122").Value = "IsMoAEInPaVILP1"
$ws.Range("A123").Value = "IsMoAEGrEqPaLP"
$ws.Range("A124").Value = "BeTaAEInPaVILP3"
$ws.Range("A125").Value = "BeTaAEInPaVILP2"
$ws.Range("A126").Value = "BeTaAEMiMaDeIIGPLP"
$ws.Range("A127").Value = "BeTaAEInSBIVPaLP"
$ws.Range("A128").Value = "BeTaAEMiMaDeIVGPLP"
$ws.Range("A129").Value = "BeTaAEInPaVLP"
$ws.Range("A130").Value = "BeTaAEInPaEFIILP"
$ws.Range("A131").Value = "BeTaAEInSBIIPaLP2"
$ws.Range("A132").Value = "BeTaAEMePaIVLP"
$ws.Range("A133").Value = "BeTaAEPaEXCFLP"
$ws.Range("A134").Value = "BeTaAEImPaLP"
$ws.Range("A135").Value = "BeTaAEInPaVILP1"
$ws.Range("A136").Value = "BeTaAEGrEqPaLP"
$ws.Range("A137").Value = "BeTaAEMePaIILP1"
$ws.Range("A138").Value = "AvAnAEInPaVILP3"
$ws.Range("A139").Value = "AvAnAEPaEXCFLP"
$ws.Range("A140").Value = "AvAnAEGrEqPaLP"
$ws.Range("A141").Value = "AvAnAEInPaVILP2"
